# Scheduled-runner style refresh of cached marketboard price/profit figures
# across the per-job Leve profit tables (columns H..N = currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 1618.75
$ws.Cells.Item(64, 9).Value = 1737.5
$ws.Cells.Item(64, 10).Value = 1500
$ws.Cells.Item(64, 11).Value = 1737.5
$ws.Cells.Item(64, 12).Value = 1500
$ws.Cells.Item(64, 13).Value = -1489.5
$ws.Cells.Item(64, 14).Value = -1996
$ws.Cells.Item(67, 8).Value = 1618.75
$ws.Cells.Item(67, 9).Value = 1737.5
$ws.Cells.Item(67, 10).Value = 1500
$ws.Cells.Item(67, 11).Value = 1737.5
$ws.Cells.Item(67, 12).Value = 1500
$ws.Cells.Item(67, 13).Value = -879.5
$ws.Cells.Item(67, 14).Value = -3216
$ws.Cells.Item(105, 8).Value = 75335.5
$ws.Cells.Item(105, 10).Value = 75335.5
$ws.Cells.Item(105, 12).Value = 75335.5
$ws.Cells.Item(105, 14).Value = -82323.5
$ws.Cells.Item(116, 8).Value = 6371.3335
$ws.Cells.Item(116, 9).Value = 5823.75
$ws.Cells.Item(116, 10).Value = 7466.5
$ws.Cells.Item(116, 11).Value = 5823.75
$ws.Cells.Item(116, 12).Value = 7466.5
$ws.Cells.Item(116, 13).Value = -2381.75
$ws.Cells.Item(116, 14).Value = -14350.5
$ws.Cells.Item(125, 8).Value = 9981.125
$ws.Cells.Item(125, 9).Value = 8712.5
$ws.Cells.Item(125, 11).Value = 78412.5
$ws.Cells.Item(125, 13).Value = -75952.5
$ws.Cells.Item(137, 8).Value = 2443.2
$ws.Cells.Item(137, 9).Value = 1738.8889
$ws.Cells.Item(137, 10).Value = 3499.6667
$ws.Cells.Item(137, 11).Value = 5216.6667
$ws.Cells.Item(137, 12).Value = 10499.0001
$ws.Cells.Item(137, 13).Value = -2666.6667
$ws.Cells.Item(137, 14).Value = -15599.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6302.7144
$ws.Cells.Item(32, 9).Value = 5217.85
$ws.Cells.Item(32, 10).Value = 28000
$ws.Cells.Item(32, 11).Value = 5217.85
$ws.Cells.Item(32, 12).Value = 28000
$ws.Cells.Item(32, 13).Value = -4930.85
$ws.Cells.Item(32, 14).Value = -28574
$ws.Cells.Item(45, 8).Value = 1528.1666
$ws.Cells.Item(45, 9).Value = 1333.8
$ws.Cells.Item(45, 11).Value = 1333.8
$ws.Cells.Item(45, 13).Value = -956.8
$ws.Cells.Item(98, 8).Value = 24999.5
$ws.Cells.Item(98, 10).Value = 24999.5
$ws.Cells.Item(98, 12).Value = 24999.5
$ws.Cells.Item(98, 14).Value = -30989.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 246.83333
$ws.Cells.Item(22, 10).Value = 167
$ws.Cells.Item(22, 12).Value = 167
$ws.Cells.Item(22, 14).Value = -513
$ws.Cells.Item(86, 8).Value = 4058.8
$ws.Cells.Item(86, 9).Value = 3948.5
$ws.Cells.Item(86, 11).Value = 3948.5
$ws.Cells.Item(86, 13).Value = -2825.5
$ws.Cells.Item(89, 8).Value = 4058.8
$ws.Cells.Item(89, 9).Value = 3948.5
$ws.Cells.Item(89, 11).Value = 19742.5
$ws.Cells.Item(89, 13).Value = -14126.5
$ws.Cells.Item(100, 8).Value = 29106
$ws.Cells.Item(100, 10).Value = 29106
$ws.Cells.Item(100, 12).Value = 29106
$ws.Cells.Item(100, 14).Value = -31270

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4249.5
$ws.Cells.Item(31, 9).Value = 5999
$ws.Cells.Item(31, 11).Value = 5999
$ws.Cells.Item(31, 13).Value = -5704
$ws.Cells.Item(34, 8).Value = 4249.5
$ws.Cells.Item(34, 9).Value = 5999
$ws.Cells.Item(34, 11).Value = 5999
$ws.Cells.Item(34, 13).Value = -5797
$ws.Cells.Item(74, 8).Value = 37754.465
$ws.Cells.Item(74, 10).Value = 37754.465
$ws.Cells.Item(74, 12).Value = 37754.465
$ws.Cells.Item(74, 14).Value = -39502.465
$ws.Cells.Item(77, 8).Value = 37754.465
$ws.Cells.Item(77, 10).Value = 37754.465
$ws.Cells.Item(77, 12).Value = 113263.395
$ws.Cells.Item(77, 14).Value = -121999.395
$ws.Cells.Item(99, 8).Value = 2499.6667
$ws.Cells.Item(99, 9).Value = 1500
$ws.Cells.Item(99, 10).Value = 2999.5
$ws.Cells.Item(99, 11).Value = 1500
$ws.Cells.Item(99, 12).Value = 2999.5
$ws.Cells.Item(99, 13).Value = -2
$ws.Cells.Item(99, 14).Value = -5995.5
$ws.Cells.Item(126, 8).Value = 2499.6667
$ws.Cells.Item(126, 9).Value = 1500
$ws.Cells.Item(126, 10).Value = 2999.5
$ws.Cells.Item(126, 11).Value = 4500
$ws.Cells.Item(126, 12).Value = 8998.5
$ws.Cells.Item(126, 13).Value = -2030
$ws.Cells.Item(126, 14).Value = -13938.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 2155.4736
$ws.Cells.Item(4, 9).Value = 1950.4445
$ws.Cells.Item(4, 10).Value = 2340
$ws.Cells.Item(4, 11).Value = 5851.333500000001
$ws.Cells.Item(4, 12).Value = 7020
$ws.Cells.Item(4, 13).Value = -5739.333500000001
$ws.Cells.Item(4, 14).Value = -7244

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(19, 8).Value = 933
$ws.Cells.Item(19, 10).Value = 1200
$ws.Cells.Item(19, 12).Value = 1200
$ws.Cells.Item(19, 14).Value = -1776
$ws.Cells.Item(132, 8).Value = 2336.3333
$ws.Cells.Item(132, 9).Value = 2403.6
$ws.Cells.Item(132, 10).Value = 2000
$ws.Cells.Item(132, 11).Value = 7210.799999999999
$ws.Cells.Item(132, 12).Value = 6000
$ws.Cells.Item(132, 13).Value = -4680.799999999999
$ws.Cells.Item(132, 14).Value = -11060

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(4, 8).Value = 1000
$ws.Cells.Item(4, 10).Value = 1000
$ws.Cells.Item(4, 12).Value = 1000
$ws.Cells.Item(4, 14).Value = -1226
$ws.Cells.Item(22, 8).Value = 960
$ws.Cells.Item(22, 10).Value = 1200
$ws.Cells.Item(22, 12).Value = 1200
$ws.Cells.Item(22, 14).Value = -1790
$ws.Cells.Item(27, 8).Value = 960
$ws.Cells.Item(27, 10).Value = 1200
$ws.Cells.Item(27, 12).Value = 1200
$ws.Cells.Item(27, 14).Value = -1414
$ws.Cells.Item(28, 8).Value = 1000
$ws.Cells.Item(28, 10).Value = 1000
$ws.Cells.Item(28, 12).Value = 1000
$ws.Cells.Item(28, 14).Value = -1464
$ws.Cells.Item(37, 8).Value = 1000
$ws.Cells.Item(37, 10).Value = 1000
$ws.Cells.Item(37, 12).Value = 1000
$ws.Cells.Item(37, 14).Value = -1214
$ws.Cells.Item(46, 8).Value = 2959.9
$ws.Cells.Item(46, 9).Value = 533.6667
$ws.Cells.Item(46, 10).Value = 3999.7144
$ws.Cells.Item(46, 11).Value = 533.6667
$ws.Cells.Item(46, 12).Value = 3999.7144
$ws.Cells.Item(46, 13).Value = -345.6667
$ws.Cells.Item(46, 14).Value = -4375.7144
$ws.Cells.Item(68, 8).Value = 5799.6665
$ws.Cells.Item(68, 9).Value = 5949.5
$ws.Cells.Item(68, 11).Value = 5949.5
$ws.Cells.Item(68, 13).Value = -5200.5
$ws.Cells.Item(71, 8).Value = 5799.6665
$ws.Cells.Item(71, 9).Value = 5949.5
$ws.Cells.Item(71, 11).Value = 29747.5
$ws.Cells.Item(71, 13).Value = -26003.5
$ws.Cells.Item(82, 8).Value = 1668.3334
$ws.Cells.Item(82, 9).Value = 1536.6666
$ws.Cells.Item(82, 10).Value = 1800
$ws.Cells.Item(82, 11).Value = 1536.6666
$ws.Cells.Item(82, 12).Value = 1800
$ws.Cells.Item(82, 13).Value = -1175.6666
$ws.Cells.Item(82, 14).Value = -2522
$ws.Cells.Item(85, 8).Value = 1668.3334
$ws.Cells.Item(85, 9).Value = 1536.6666
$ws.Cells.Item(85, 10).Value = 1800
$ws.Cells.Item(85, 11).Value = 1536.6666
$ws.Cells.Item(85, 12).Value = 1800
$ws.Cells.Item(85, 13).Value = -288.6666
$ws.Cells.Item(85, 14).Value = -4296

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1179
$ws.Cells.Item(132, 9).Value = 1179
$ws.Cells.Item(132, 11).Value = 3537
$ws.Cells.Item(132, 13).Value = -1007
